$d = $word.ActiveDocument

# --- Locate the target sentence precisely (avoid collisions with similar text elsewhere) ---
$anchor = $d.Content
$anchor.Find.Execute("շարժիչը (Rotate_Engine()) 8-րդ", $true, $false, $false, $false, $false,
                      $true, 1, $false, "", 0)
$sentenceStart = $anchor.Start

# 1) "Պտտում" -> "Տեղաշարժում"  (simple in-run text replace)
$d.Content.Find.Execute("Պտտում", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Տեղաշարժում", 2)

# 2) "շարժիչը" -> "սանրը"  (simple in-run text replace)
$d.Content.Find.Execute("շարժիչը", $true, $false, $false, $false, $false,
                         $true, 1, $false, "սանրը", 2)

# 3) Split the " (" run into " " and "(" runs (toggling bold on/off forces the
#    run-boundary split without leaving any visible formatting artifact).
$r = $d.Content
$r.Find.Execute("սանրը (")
$parenStart = $r.End - 1
$parenRange = $d.Range($parenStart, $parenStart + 1)
$parenRange.Font.Bold = $true
$parenRange.Font.Bold = $false

# 4) "Rotate_" -> "Move_"  (simple in-run text replace)
$d.Content.Find.Execute("Rotate_", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Move_", 2)

# 5) "Engine" -> "Brush"  (simple in-run text replace)
$d.Content.Find.Execute("Engine", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Brush", 2)

# 6) Split the ")) 8-րդ " run into ")" and ") 8-րդ " runs the same way.
$r2 = $d.Content
$r2.Find.Execute("Brush(")
$closeStart = $r2.End
$closeRange = $d.Range($closeStart, $closeStart + 1)
$closeRange.Font.Bold = $true
$closeRange.Font.Bold = $false
